# Generate Report for Handoff
# - Flip the "Status" text from "Handed back: in sync with en-US" to
#   "Ready for handoff" everywhere it appears (Overview + per-locale sheets).
# - Bump the associated timestamps to reflect the new handoff generation.
# - Column widths for the Status columns auto-shrink since the new text is
#   shorter than the old text.

$wb = $excel.ActiveWorkbook

$ovw  = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$ovw.Range("E2").Value = "Ready for handoff"
$ovw.Range("F2").Value = "Ready for handoff"
$ovw.Range("G2").Value = "2016-08-31 17:10:32"

# --- zh-cn detail sheet ------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-31 17:10:27"

# --- de-de detail sheet ------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-31 17:10:32"

# --- Column widths: Status columns shrink now that the text is shorter -
$ovw.Columns.Item(5).ColumnWidth = 16.38265482584637
$ovw.Columns.Item(6).ColumnWidth = 16.38265482584637
$zhcn.Columns.Item(3).ColumnWidth = 16.38265482584637
$dede.Columns.Item(3).ColumnWidth = 16.38265482584637
